$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Try to avoid conflicts if / you are / applying to the same department
#    (there are plenty of universities with essentially the same external
#    profile). " -- merge the three runs into a single run by replacing the
#    combined text with itself (Word collapses same-format adjacent runs
#    touched by a single Find/Replace into one run).
# ---------------------------------------------------------------------------
$tryAvoidText = "Try to avoid conflicts if you are applying to the same department (there are plenty of universities with essentially the same external profile). "
$d.Content.Find.Execute($tryAvoidText, $true, $false, $false, $false, $false, $true, 1, $false, $tryAvoidText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "In Europe, / it will / be mostly project specific, so your boss will
#    likely be fixed in advance (there are exceptions of course). " -- same
#    merge treatment.
# ---------------------------------------------------------------------------
$europeText = "In Europe, it will be mostly project specific, so your boss will likely be fixed in advance (there are exceptions of course). "
$d.Content.Find.Execute($europeText, $true, $false, $false, $false, $false, $true, 1, $false, $europeText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "For some people, the location matters. If / you are / into hiking, or
#    outdoor sport, or if you enjoy specific activities which you feel are
#    essential to your happiness (" -- same merge treatment (stop right
#    before the proofErr-wrapped "I'm" run, which must stay untouched). This
#    paragraph's leading run also carries a <w:lastRenderedPageBreak/>
#    marker that a plain Find/Replace across the run boundaries would drop,
#    so the whole paragraph is rebuilt explicitly via InsertXML instead,
#    leaving the proofErr/"I'm"/trailing-sentence runs untouched.
# ---------------------------------------------------------------------------
$locationText = "For some people, the location matters. If you are into hiking, or outdoor sport, or if you enjoy specific activities which you feel are essential to your happiness ("

$rng2 = $d.Content
$rng2.Find.Execute("factor that in.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$locPara = $rng2.Paragraphs(1)
$locRange = $locPara.Range

$wNs2 = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$locPPr = "<w:pPr><w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:ind w:left='540'/><w:textAlignment w:val='center'/><w:rPr><w:rFonts w:ascii='Calibri' w:eastAsia='Times New Roman' w:hAnsi='Calibri' w:cs='Calibri'/><w:szCs w:val='22'/><w:lang w:eastAsia='en-GB'/></w:rPr></w:pPr>"
$locRunRPr = "<w:rPr><w:rFonts w:ascii='Calibri' w:eastAsia='Times New Roman' w:hAnsi='Calibri' w:cs='Calibri'/><w:szCs w:val='22'/><w:lang w:eastAsia='en-GB'/></w:rPr>"

$locRun1 = "<w:r>$locRunRPr<w:lastRenderedPageBreak/><w:t>$locationText</w:t></w:r>"
$locProof1 = "<w:proofErr w:type='gramStart'/>"
$locRun2 = "<w:r>$locRunRPr<w:t>I'm</w:t></w:r>"
$locProof2 = "<w:proofErr w:type='gramEnd'/>"
$locRun3 = "<w:r>$locRunRPr<w:t xml:space='preserve'> a robot, so I had no requirements), factor that in.</w:t></w:r>"

$locXml = "<w:p $wNs2>$locPPr$locRun1$locProof1$locRun2$locProof2$locRun3</w:p>"
$locRange.InsertXML($locXml) | Out-Null

# ---------------------------------------------------------------------------
# 4) The paragraph that used to read "Avoid formatting errors, ... targeting
#    the right level." (three runs) now becomes two paragraphs:
#      a) "Ask your seniors/alumni for places if you know your field of
#         interest; they would usually know where the good groups are
#         (because they are actively working in the field)." -- with a
#         slightly trimmed run-level rPr (no eastAsia font / no eastAsia
#         lang override).
#      b) a new paragraph carrying the original "Avoid formatting errors...
#         targeting the right level." sentence, merged into a single run,
#         right before "Remember: Professors are looking..."
# ---------------------------------------------------------------------------
$avoidText = "Avoid formatting errors, spelling errors, especially in names. Try to re-read the application after a few hours before finally submitting it. Ensure your SOPs are peer-reviewed. Perhaps have one quickly glanced at by a professor to ensure you are targeting the right level."
$askText = "Ask your seniors/alumni for places if you know your field of interest; they would usually know where the good groups are (because they are actively working in the field)."

$rng = $d.Content
$rng.Find.Execute($avoidText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)
$r = $para.Range

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$pPrCommon = "<w:pPr><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr><w:spacing w:after='0' w:line='240' w:lineRule='auto'/><w:ind w:left='540'/><w:textAlignment w:val='center'/><w:rPr><w:rFonts w:ascii='Calibri' w:eastAsia='Times New Roman' w:hAnsi='Calibri' w:cs='Calibri'/><w:szCs w:val='22'/><w:lang w:eastAsia='en-GB'/></w:rPr></w:pPr>"

$askRunRPr = "<w:rPr><w:rFonts w:ascii='Calibri' w:hAnsi='Calibri' w:cs='Calibri'/><w:szCs w:val='22'/></w:rPr>"
$avoidRunRPr = "<w:rPr><w:rFonts w:ascii='Calibri' w:eastAsia='Times New Roman' w:hAnsi='Calibri' w:cs='Calibri'/><w:szCs w:val='22'/><w:lang w:eastAsia='en-GB'/></w:rPr>"

$p1 = "<w:p $wNs>$pPrCommon<w:r>$askRunRPr<w:t>$askText</w:t></w:r></w:p>"
$p2 = "<w:p $wNs>$pPrCommon<w:r>$avoidRunRPr<w:t>$avoidText</w:t></w:r></w:p>"

$r.InsertXML($p1 + $p2) | Out-Null

Write-Output "done"
